$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 470, shifting all rows from
# 470 downward (through 565) down to 472..567.
$ws.Rows("470:471").Insert()

# New row 470 (new weekly observation)
$ws.Cells.Item(470, 1).Value = 9
$ws.Cells.Item(470, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(470, 3).Value = "Metropolitana"
$ws.Cells.Item(470, 4).Value = 44782
$ws.Cells.Item(470, 5).Value = 13
$ws.Cells.Item(470, 6).Value = 100112031
$ws.Cells.Item(470, 7).Value = "Poroto verde"
$ws.Cells.Item(470, 8).Value = "Magnum"
$ws.Cells.Item(470, 9).Value = "Primera"
$ws.Cells.Item(470, 10).Value = 67
$ws.Cells.Item(470, 11).Value = 36000
$ws.Cells.Item(470, 12).Value = 38000
$ws.Cells.Item(470, 13).Value = 37045
$ws.Cells.Item(470, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(470, 15).Value = "Perú"
$ws.Cells.Item(470, 16).Value = 1482
$ws.Cells.Item(470, 17).Value = 25
$ws.Cells.Item(470, 18).Value = "Hortaliza"

# New row 471 (new weekly observation)
$ws.Cells.Item(471, 1).Value = 9
$ws.Cells.Item(471, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(471, 3).Value = "Metropolitana"
$ws.Cells.Item(471, 4).Value = 44782
$ws.Cells.Item(471, 5).Value = 13
$ws.Cells.Item(471, 6).Value = 100112031
$ws.Cells.Item(471, 7).Value = "Poroto verde"
$ws.Cells.Item(471, 8).Value = "Sin especificar"
$ws.Cells.Item(471, 9).Value = "Primera"
$ws.Cells.Item(471, 10).Value = 52
$ws.Cells.Item(471, 11).Value = 44000
$ws.Cells.Item(471, 12).Value = 45000
$ws.Cells.Item(471, 13).Value = 44385
$ws.Cells.Item(471, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(471, 15).Value = "Perú"
$ws.Cells.Item(471, 16).Value = 1775
$ws.Cells.Item(471, 17).Value = 25
$ws.Cells.Item(471, 18).Value = "Hortaliza"
